# Generate Report for Handoff
# Updates the status of the localization from "In Translation" to
# "Ready for handoff" and refreshes the related timestamps on the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 08:47:27"

# zh-cn sheet: Status column (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 08:47:22"

# de-de sheet: Status column (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 08:47:27"

# Widen the status columns to fit the new, longer text ("Ready for
# handoff" vs. "In Translation"). 16.3 is the COM-level ColumnWidth that
# maps to the closest achievable underlying column width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
